# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# Rows: 2,3,4,6,7,8 change; row 5 stays the same.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1345
    $ws.Range("F3").Value = 1911
    $ws.Range("F4").Value = 190
    $ws.Range("F6").Value = 6324
    $ws.Range("F7").Value = 192
    $ws.Range("F8").Value = 113
}
